# Add a new attendance sheet "2021-09-18" after the last existing
# attendance sheet ("2021-09-10"), reusing its layout/formatting, then
# fill in the new day's values.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-09-10")

# Copy the template sheet so the new sheet inherits the same column
# headers, cell styles (bold/centered/bordered header row + A2) and page
# setup, and place the copy right after the template (i.e. at the end).
$template.Copy($null, $template)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2021-09-18"

# Update the day-specific values in row 2; everything else (headers,
# Face_ID, Name, Address, Job, Heart-rate-less columns) stays identical
# to the template sheet.
$ws.Range("E2").Value = "17:10:59"
$ws.Range("F2").Value = 98.83116074550881
$ws.Range("G2").Value = 78

# Copying a sheet activates it; restore the original active tab (Sheet1)
# so we only add the new sheet without disturbing the saved selection.
$wb.Worksheets.Item("Sheet1").Activate()
